$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '36.614.16'
$ws.Cells.Item(2, 5).Value = '  +1.41%  '
$ws.Cells.Item(3, 4).Value = '1.948.34'
$ws.Cells.Item(3, 5).Value = '  +0.87%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '244.04'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.46%  '
$ws.Cells.Item(6, 5).Value = '  +1.75%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '58.17'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +4.52%  '
$ws.Cells.Item(8, 5).Value = '  -0.04%  '
$ws.Cells.Item(9, 5).Value = '  +1.46%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0805'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -2.09%  '
$ws.Cells.Item(11, 5).Value = '  +0.14%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '22.26'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +7.21%  '
$ws.Cells.Item(13, 4).Value = '2.235.81'
$ws.Cells.Item(13, 5).Value = '  +0.81%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.815'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.47%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.49'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +2.62%  '
$ws.Cells.Item(16, 5).Value = '  +1.27%  '
$ws.Cells.Item(17, 4).Value = '1.953.11'
$ws.Cells.Item(17, 5).Value = '  +2.07%  '
$ws.Cells.Item(18, 4).Value = '36.575.18'
$ws.Cells.Item(18, 5).Value = '  +1.55%  '
$ws.Cells.Item(19, 5).Value = '  +0.27%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0857'
$ws.Cells.Item(20, 5).Value = '  -0.03%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '228.41'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +1.63%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.02'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +2.08%  '
$ws.Cells.Item(24, 5).Value = '  +0.32%  '
$ws.Cells.Item(25, 5).Value = '  +3.29%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.23'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -0.02%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '160.24'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -1.32%  '
$ws.Cells.Item(28, 5).Value = '  +15.38%  '
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.34'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +1.22%  '
$ws.Cells.Item(30, 5).Value = '  +1.85%  '
$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.68'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.80%  '
$ws.Cells.Item(32, 2).Value = 'ImmutableX'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.10'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -2.05%  '
$ws.Cells.Item(33, 5).Value = '  +0.50%  '
$ws.Cells.Item(34, 5).Value = '  -0.66%  '
$ws.Cells.Item(35, 5).Value = '  +6.17%  '
$ws.Cells.Item(36, 5).Value = '  +0.18%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.44'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +21.93%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.21'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +3.88%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.77'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.75%  '
$ws.Cells.Item(40, 5).Value = '  +3.42%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.90'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +1.12%  '
$ws.Cells.Item(42, 5).Value = '  +2.61%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.16'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.38%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '15.94'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +3.49%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.03'
$cell.Style = "Normal"
$ws.Cells.Item(46, 4).Value = '1.343.24'
$ws.Cells.Item(46, 5).Value = '  +0.72%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '86.99'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.01%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.17'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +0.19%  '
$ws.Cells.Item(49, 5).Value = '  +1.05%  '
$ws.Cells.Item(50, 4).Value = '2.127.65'
$ws.Cells.Item(50, 5).Value = '  +0.80%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '43.29'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -5.37%  '
